$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 612.5
$ws.Range("I33").Value = 213.63637
$ws.Range("K33").Value = 213.63637
$ws.Range("M33").Value = 15.36363

$ws.Range("H64").Value = 23029.154
$ws.Range("I64").Value = 3772.261
$ws.Range("K64").Value = 3772.261
$ws.Range("M64").Value = -3524.261

$ws.Range("H67").Value = 23029.154
$ws.Range("I67").Value = 3772.261
$ws.Range("K67").Value = 3772.261
$ws.Range("M67").Value = -2914.261

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H98").Value = 1351.7391
$ws.Range("I98").Value = 1331.3636
$ws.Range("J98").Value = 1800
$ws.Range("K98").Value = 1331.3636
$ws.Range("L98").Value = 1800
$ws.Range("M98").Value = 166.6364000000001
$ws.Range("N98").Value = -4796

$ws.Range("H100").Value = 12144.429
$ws.Range("I100").Value = 3005
$ws.Range("K100").Value = 3005
$ws.Range("M100").Value = -2464

$ws.Range("H122").Value = 1351.7391
$ws.Range("I122").Value = 1331.3636
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 3994.0908
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -1544.0908
$ws.Range("N122").Value = -10300

$ws.Range("H125").Value = 71429896
$ws.Range("I125").Value = 83334710
$ws.Range("J125").Value = 1000
$ws.Range("K125").Value = 750012390
$ws.Range("L125").Value = 9000
$ws.Range("M125").Value = -750009930
$ws.Range("N125").Value = -13920

$ws.Range("H141").Value = 2866.5386
$ws.Range("I141").Value = 2473.2144
$ws.Range("J141").Value = 3325.4167
$ws.Range("K141").Value = 7419.6432
$ws.Range("L141").Value = 9976.250100000001
$ws.Range("M141").Value = -2239.6432
$ws.Range("N141").Value = -20336.2501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1340
$ws.Range("I2").Value = 900
$ws.Range("K2").Value = 900
$ws.Range("M2").Value = -787

$ws.Range("H5").Value = 1881.2667
$ws.Range("I5").Value = 1607.2727
$ws.Range("J5").Value = 2634.75
$ws.Range("K5").Value = 1607.2727
$ws.Range("L5").Value = 2634.75
$ws.Range("M5").Value = -1495.2727
$ws.Range("N5").Value = -2858.75

$ws.Range("H16").Value = 28547.166
$ws.Range("I16").Value = 425.33334
$ws.Range("J16").Value = 56669
$ws.Range("K16").Value = 425.33334
$ws.Range("L16").Value = 56669
$ws.Range("M16").Value = -138.33334
$ws.Range("N16").Value = -57243

$ws.Range("H32").Value = 212618.36
$ws.Range("I32").Value = 203354.92
$ws.Range("J32").Value = 367009
$ws.Range("K32").Value = 203354.92
$ws.Range("L32").Value = 367009
$ws.Range("M32").Value = -203067.92
$ws.Range("N32").Value = -367583

$ws.Range("H61").Value = 2243.3076
$ws.Range("J61").Value = 1457
$ws.Range("L61").Value = 1457
$ws.Range("N61").Value = -1881

$ws.Range("H116").Value = 1340
$ws.Range("I116").Value = 900
$ws.Range("K116").Value = 900
$ws.Range("M116").Value = 1394

$ws.Range("H136").Value = 2243.3076
$ws.Range("J136").Value = 1457
$ws.Range("L136").Value = 4371
$ws.Range("N136").Value = -9471

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1340
$ws.Range("I3").Value = 900
$ws.Range("K3").Value = 900
$ws.Range("M3").Value = -786

$ws.Range("H4").Value = 1881.2667
$ws.Range("I4").Value = 1607.2727
$ws.Range("J4").Value = 2634.75
$ws.Range("K4").Value = 1607.2727
$ws.Range("L4").Value = 2634.75
$ws.Range("M4").Value = -1492.2727
$ws.Range("N4").Value = -2864.75

$ws.Range("H20").Value = 1783.3334
$ws.Range("I20").Value = 1775
$ws.Range("J20").Value = 1800
$ws.Range("K20").Value = 1775
$ws.Range("L20").Value = 1800
$ws.Range("M20").Value = -1528
$ws.Range("N20").Value = -2294

$ws.Range("H80").Value = 299.36365
$ws.Range("J80").Value = 304.25
$ws.Range("L80").Value = 304.25
$ws.Range("N80").Value = -2300.25

$ws.Range("H83").Value = 299.36365
$ws.Range("J83").Value = 304.25
$ws.Range("L83").Value = 1521.25
$ws.Range("N83").Value = -11505.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1904.0264
$ws.Range("I132").Value = 1211.9656
$ws.Range("J132").Value = 4134
$ws.Range("K132").Value = 3635.8968
$ws.Range("L132").Value = 12402
$ws.Range("M132").Value = -1105.8968
$ws.Range("N132").Value = -17462

$ws.Range("H134").Value = 14707536
$ws.Range("I134").Value = 1415.3928
$ws.Range("J134").Value = 83336100
$ws.Range("K134").Value = 4246.178400000001
$ws.Range("L134").Value = 250008300
$ws.Range("M134").Value = -1711.178400000001
$ws.Range("N134").Value = -250013370

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 4843.5557
$ws.Range("J104").Value = 4843.5557
$ws.Range("L104").Value = 14530.6671
$ws.Range("N104").Value = -19772.6671

$ws.Range("H122").Value = 444.0909
$ws.Range("I122").Value = 310.6111
$ws.Range("J122").Value = 1044.75
$ws.Range("K122").Value = 2795.4999
$ws.Range("L122").Value = 9402.75
$ws.Range("M122").Value = -345.4999000000003
$ws.Range("N122").Value = -14302.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 21845
$ws.Range("I31").Value = 2976
$ws.Range("J31").Value = 69017.5
$ws.Range("K31").Value = 2976
$ws.Range("L31").Value = 69017.5
$ws.Range("M31").Value = -2684
$ws.Range("N31").Value = -69601.5

$ws.Range("H37").Value = 21845
$ws.Range("I37").Value = 2976
$ws.Range("J37").Value = 69017.5
$ws.Range("K37").Value = 2976
$ws.Range("L37").Value = 69017.5
$ws.Range("M37").Value = -2699
$ws.Range("N37").Value = -69571.5

$ws.Range("H102").Value = 24696.154
$ws.Range("I102").Value = 1737.3334
$ws.Range("J102").Value = 76353.5
$ws.Range("K102").Value = 1737.3334
$ws.Range("L102").Value = 76353.5
$ws.Range("M102").Value = -115.3334
$ws.Range("N102").Value = -79597.5

$ws.Range("H113").Value = 2790.625
$ws.Range("I113").Value = 2825
$ws.Range("J113").Value = 2756.25
$ws.Range("K113").Value = 2825
$ws.Range("L113").Value = 2756.25
$ws.Range("M113").Value = -655
$ws.Range("N113").Value = -7096.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2993.0344
$ws.Range("I7").Value = 3098.15
$ws.Range("J7").Value = 2759.4443
$ws.Range("K7").Value = 3098.15
$ws.Range("L7").Value = 2759.4443
$ws.Range("M7").Value = -2986.15
$ws.Range("N7").Value = -2983.4443

$ws.Range("H64").Value = 20000
$ws.Range("J64").Value = 20000
$ws.Range("L64").Value = 20000
$ws.Range("N64").Value = -20450

$ws.Range("H67").Value = 20000
$ws.Range("J67").Value = 20000
$ws.Range("L67").Value = 20000
$ws.Range("N67").Value = -21560

$ws.Range("H126").Value = 2993.0344
$ws.Range("I126").Value = 3098.15
$ws.Range("J126").Value = 2759.4443
$ws.Range("K126").Value = 9294.450000000001
$ws.Range("L126").Value = 8278.332900000001
$ws.Range("M126").Value = -6824.450000000001
$ws.Range("N126").Value = -13218.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 69336
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 69336
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 69336
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -69616

$ws.Range("H126").Value = 1461.8462
$ws.Range("I126").Value = 952
$ws.Range("J126").Value = 1554.5454
$ws.Range("K126").Value = 2856
$ws.Range("L126").Value = 4663.6362
$ws.Range("M126").Value = -386
$ws.Range("N126").Value = -9603.636200000001

